$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 329.7143
$ws.Range("I5").Value = 45.2
$ws.Range("K5").Value = 45.2
$ws.Range("M5").Value = 69.8

$ws.Range("H17").Value = 1176.6842
$ws.Range("J17").Value = 1176.6842
$ws.Range("L17").Value = 3530.0526
$ws.Range("N17").Value = -3866.0526

$ws.Range("H64").Value = 37379.586
$ws.Range("I64").Value = 93551
$ws.Range("K64").Value = 93551
$ws.Range("M64").Value = -93303

$ws.Range("H67").Value = 37379.586
$ws.Range("I67").Value = 93551
$ws.Range("K67").Value = 93551
$ws.Range("M67").Value = -92693

$ws.Range("H86").Value = 30315.943
$ws.Range("I86").Value = 47319.863
$ws.Range("J86").Value = 1540.0769
$ws.Range("K86").Value = 47319.863
$ws.Range("L86").Value = 1540.0769
$ws.Range("M86").Value = -46196.863
$ws.Range("N86").Value = -3786.0769

$ws.Range("H87").Value = 31878.5
$ws.Range("J87").Value = 31878.5
$ws.Range("L87").Value = 31878.5
$ws.Range("N87").Value = -34374.5

$ws.Range("H89").Value = 30315.943
$ws.Range("I89").Value = 47319.863
$ws.Range("J89").Value = 1540.0769
$ws.Range("K89").Value = 236599.315
$ws.Range("L89").Value = 7700.3845
$ws.Range("M89").Value = -230983.315
$ws.Range("N89").Value = -18932.3845

$ws.Range("H90").Value = 31878.5
$ws.Range("J90").Value = 31878.5
$ws.Range("L90").Value = 95635.5
$ws.Range("N90").Value = -108115.5

$ws.Range("H103").Value = 1220.7273
$ws.Range("I103").Value = 2787.5
$ws.Range("K103").Value = 8362.5
$ws.Range("M103").Value = -7776.5

$ws.Range("H113").Value = 168800.67
$ws.Range("J113").Value = 1899.5
$ws.Range("L113").Value = 1899.5
$ws.Range("N113").Value = -8407.5

$ws.Range("H132").Value = 4102508
$ws.Range("I132").Value = 5004424.5
$ws.Range("K132").Value = 15013273.5
$ws.Range("M132").Value = -15010743.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2542.2222
$ws.Range("I63").Value = 1980
$ws.Range("J63").Value = 3666.6667
$ws.Range("K63").Value = 1980
$ws.Range("L63").Value = 3666.6667
$ws.Range("M63").Value = -1294
$ws.Range("N63").Value = -5038.6667

$ws.Range("H66").Value = 2542.2222
$ws.Range("I66").Value = 1980
$ws.Range("J66").Value = 3666.6667
$ws.Range("K66").Value = 9900
$ws.Range("L66").Value = 18333.3335
$ws.Range("M66").Value = -6468
$ws.Range("N66").Value = -25197.3335

$ws.Range("H74").Value = 965.6842
$ws.Range("I74").Value = 878.1667
$ws.Range("K74").Value = 878.1667
$ws.Range("M74").Value = -4.166699999999992

$ws.Range("H77").Value = 965.6842
$ws.Range("I77").Value = 878.1667
$ws.Range("K77").Value = 4390.8335
$ws.Range("M77").Value = -22.83349999999973

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 200
$ws.Range("I5").Value = 200
$ws.Range("K5").Value = 200
$ws.Range("M5").Value = -87

$ws.Range("H134").Value = 2086.853
$ws.Range("I134").Value = 1713.2667
$ws.Range("J134").Value = 4888.75
$ws.Range("K134").Value = 5139.800099999999
$ws.Range("L134").Value = 14666.25
$ws.Range("M134").Value = -2604.800099999999
$ws.Range("N134").Value = -19736.25

$ws.Range("H140").Value = 38819.9
$ws.Range("J140").Value = 40832.223
$ws.Range("L140").Value = 40832.223
$ws.Range("N140").Value = -51192.223

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 2402900
$ws.Range("J4").Value = 503625
$ws.Range("L4").Value = 503625
$ws.Range("N4").Value = -503849

$ws.Range("H7").Value = 116.76471
$ws.Range("I7").Value = 67.90000000000001
$ws.Range("J7").Value = 186.57143
$ws.Range("K7").Value = 67.90000000000001
$ws.Range("L7").Value = 186.57143
$ws.Range("M7").Value = 45.09999999999999
$ws.Range("N7").Value = -412.57143

$ws.Range("H140").Value = 54839.4
$ws.Range("J140").Value = 54839.4
$ws.Range("L140").Value = 54839.4
$ws.Range("N140").Value = -65199.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 852.01
$ws.Range("J131").Value = 856.57574
$ws.Range("L131").Value = 2569.72722
$ws.Range("N131").Value = -12649.72722

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 140249.6
$ws.Range("I70").Value = 205882.4
$ws.Range("J70").Value = 8984
$ws.Range("K70").Value = 205882.4
$ws.Range("L70").Value = 8984
$ws.Range("M70").Value = -205612.4
$ws.Range("N70").Value = -9524

$ws.Range("H73").Value = 140249.6
$ws.Range("I73").Value = 205882.4
$ws.Range("J73").Value = 8984
$ws.Range("K73").Value = 205882.4
$ws.Range("L73").Value = 8984
$ws.Range("M73").Value = -204946.4
$ws.Range("N73").Value = -10856

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1416.4706
$ws.Range("J22").Value = 962.8570999999999
$ws.Range("L22").Value = 962.8570999999999
$ws.Range("N22").Value = -1552.8571

$ws.Range("H27").Value = 1416.4706
$ws.Range("J27").Value = 962.8570999999999
$ws.Range("L27").Value = 962.8570999999999
$ws.Range("N27").Value = -1176.8571

$ws.Range("H68").Value = 4820.8335
$ws.Range("I68").Value = 2366.6667
$ws.Range("J68").Value = 5638.8887
$ws.Range("K68").Value = 2366.6667
$ws.Range("L68").Value = 5638.8887
$ws.Range("M68").Value = -1617.6667
$ws.Range("N68").Value = -7136.8887

$ws.Range("H71").Value = 4820.8335
$ws.Range("I71").Value = 2366.6667
$ws.Range("J71").Value = 5638.8887
$ws.Range("K71").Value = 11833.3335
$ws.Range("L71").Value = 28194.4435
$ws.Range("M71").Value = -8089.333500000001
$ws.Range("N71").Value = -35682.4435

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9617659
$ws.Range("J62").Value = 2600
$ws.Range("L62").Value = 2600
$ws.Range("N62").Value = -3848

$ws.Range("H65").Value = 9617659
$ws.Range("J65").Value = 2600
$ws.Range("L65").Value = 13000
$ws.Range("N65").Value = -19240

$ws.Range("H132").Value = 2172.2856
$ws.Range("I132").Value = 2266.6052
$ws.Range("J132").Value = 1973.1666
$ws.Range("K132").Value = 6799.8156
$ws.Range("L132").Value = 5919.4998
$ws.Range("M132").Value = -4269.8156
$ws.Range("N132").Value = -10979.4998

$ws.Range("H136").Value = 1350.0143
$ws.Range("I136").Value = 499.05554
$ws.Range("J136").Value = 2251.0293
$ws.Range("K136").Value = 1497.16662
$ws.Range("L136").Value = 6753.0879
$ws.Range("M136").Value = 1052.83338
$ws.Range("N136").Value = -11853.0879
